# Insert a new data row into the "Vega Modelo de Temuco - Zanahoria" sheet
# at position 158 (pushing the existing rows 158..249 down to 159..250),
# and populate it with the new observation dated 2022-01-11 (serial 44572).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 158; everything below shifts down one row.
$ws.Rows.Item(158).Insert()

# Fill in the values for the newly inserted row 158.
$ws.Cells.Item(158, 1).Value = 10
$ws.Cells.Item(158, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(158, 3).Value = "La Araucanía"
$ws.Cells.Item(158, 4).Value = 44572
$ws.Cells.Item(158, 5).Value = 9
$ws.Cells.Item(158, 6).Value = 100114013
$ws.Cells.Item(158, 7).Value = "Zanahoria"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 110
$ws.Cells.Item(158, 11).Value = 10000
$ws.Cells.Item(158, 12).Value = 10000
$ws.Cells.Item(158, 13).Value = 10000
$ws.Cells.Item(158, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(158, 15).Value = "Región del Maule"
$ws.Cells.Item(158, 16).Value = 500
$ws.Cells.Item(158, 17).Value = 20
$ws.Cells.Item(158, 18).Value = "Hortaliza"
